# Generate Report for handoff
# Sets the "Latest Handoff Datetime" (column D, row 3) for the
# 55d6fdfa-... file entry on both the zh-cn and de-de sheets to its own
# freshly-generated handoff timestamp (previously it just echoed the
# Datetime of the row above).

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D3").Value = "2016-01-20 07:10:20"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D3").Value = "2016-01-20 07:10:31"
